# Update column G ("K") values on Sheet1 with regenerated strikeout counts.
# These replace the previous "Strike#" derived values with true K (strikeouts)
# values as part of regenerating save_data (see commit message).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$gValues = @{
    2  = 9
    3  = 5
    4  = 4
    5  = 9
    6  = 5
    7  = 9
    8  = 6
    9  = 10
    10 = 8
    11 = 9
    12 = 5
    13 = 8
    14 = 3
    15 = 5
    16 = 8
    17 = 5
    18 = 7
    19 = 7
    20 = 8
    21 = 1
    22 = 10
    23 = 5
    24 = 4
    25 = 3
    26 = 5
    27 = 6
    28 = 4
    29 = 7
    30 = 1
    31 = 6
    32 = 2
    33 = 3
    34 = 2
    35 = 8
    36 = 1
    37 = 2
}

foreach ($row in $gValues.Keys) {
    $ws.Range("G$row").Value = $gValues[$row]
}
